$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reslice the reward-table symbols (line reward_table, reslice the images) ---
# Column A: symbol legend (rows 1-7)
$ws.Range("A3").Value = "🍎"
$ws.Range("A4").Value = "🫐"
$ws.Range("A5").Value = "🍀"
$ws.Range("A6").Value = "💰"
$ws.Range("A7").Value = "🔔"

# Row 2: triple-bell jackpot
$ws.Range("D2").Value = "🔔🔔🔔"
$ws.Range("F2").Value = "🔔🔔🔔"
$ws.Range("H2").Value = "🔔🔔🔔"
$ws.Range("J2").Value = "🔔🔔🔔"

# Row 3: triple-money-bag
$ws.Range("D3").Value = "💰💰💰"
$ws.Range("F3").Value = "💰💰💰"
$ws.Range("H3").Value = "💰💰💰"
$ws.Range("J3").Value = "💰💰💰"

# Row 4: double-money-bag + bell
$ws.Range("D4").Value = "💰💰🔔"
$ws.Range("F4").Value = "💰💰🔔"
$ws.Range("H4").Value = "💰💰🔔"
$ws.Range("J4").Value = "💰💰🔔"

# Row 5: triple-clover
$ws.Range("D5").Value = "🍀🍀🍀"
$ws.Range("F5").Value = "🍀🍀🍀"
$ws.Range("H5").Value = "🍀🍀🍀"
$ws.Range("J5").Value = "🍀🍀🍀"

# Row 6: double-clover + bell
$ws.Range("D6").Value = "🍀🍀🔔"
$ws.Range("F6").Value = "🍀🍀🔔"
$ws.Range("H6").Value = "🍀🍀🔔"
$ws.Range("J6").Value = "🍀🍀🔔"

# Row 7: triple-blueberry
$ws.Range("D7").Value = "🫐🫐🫐"
$ws.Range("F7").Value = "🫐🫐🫐"
$ws.Range("H7").Value = "🫐🫐🫐"
$ws.Range("J7").Value = "🫐🫐🫐"

# Row 8: double-blueberry + bell
$ws.Range("D8").Value = "🫐🫐🔔"
$ws.Range("F8").Value = "🫐🫐🔔"
$ws.Range("H8").Value = "🫐🫐🔔"
$ws.Range("J8").Value = "🫐🫐🔔"

# Row 9: triple-apple
$ws.Range("D9").Value = "🍎🍎🍎"
$ws.Range("F9").Value = "🍎🍎🍎"
$ws.Range("H9").Value = "🍎🍎🍎"
$ws.Range("J9").Value = "🍎🍎🍎"

# Row 10: double-apple + bell
$ws.Range("D10").Value = "🍎🍎🔔"
$ws.Range("F10").Value = "🍎🍎🔔"
$ws.Range("H10").Value = "🍎🍎🔔"
$ws.Range("J10").Value = "🍎🍎🔔"

# Row 13: cherry + ANY + ANY
$ws.Range("D13").Value = "🍒ANYANY"
$ws.Range("F13").Value = "🍒ANYANY"
$ws.Range("H13").Value = "🍒ANYANY"
$ws.Range("J13").Value = "🍒ANYANY"

# --- Move the active selection to H21 ---
$ws.Range("H21").Select()
